$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns before G (support_span, load_span)
$ws.Columns("G:H").Insert()

# Insert six new columns before the former elastic_mod_mean column (now L)
$ws.Columns("L:Q").Insert()

# --- Header row ---
$ws.Range("G1").Value = "support_span"
$ws.Range("H1").Value = "load_span"
$ws.Range("I1").Value = "char_strength"
$ws.Range("J1").Value = "weibull_modulus"
$ws.Range("L1").Value = "mean_gage_V"
$ws.Range("M1").Value = "mean_gage_SA"
$ws.Range("N1").Value = "mean_effective_V"
$ws.Range("O1").Value = "mean_effective_SA"
$ws.Range("P1").Value = "SA_series_scale_param"
$ws.Range("Q1").Value = "V_series_scale_param"

# --- Data rows ---
$ws.Range("G2").Value = 185
$ws.Range("H2").Value = 70
$ws.Range("L2").Value = 4001.743631878642
$ws.Range("M2").Value = 3050.109475517258
$ws.Range("N2").Value = 137.5265212300033
$ws.Range("O2").Value = 300.4577935658554
$ws.Range("P2").Value = 44356.98626751734
$ws.Range("Q2").Value = 20303.22442703292
$ws.Range("G3").Value = 146
$ws.Range("H3").Value = 90
$ws.Range("L3").Value = 303.7279309294768
$ws.Range("M3").Value = 746.4895383827387
$ws.Range("N3").Value = 15.59242455167278
$ws.Range("O3").Value = 104.6741842613235
$ws.Range("P3").Value = 25445.96491811102
$ws.Range("Q3").Value = 3790.469359088768
$ws.Range("G4").Value = 185
$ws.Range("H4").Value = 70
$ws.Range("L4").Value = 2280.817582212344
$ws.Range("M4").Value = 2302.69316730171
$ws.Range("N4").Value = 28.77239337819473
$ws.Range("O4").Value = 142.7393917410353
$ws.Range("P4").Value = 19154.69165355318
$ws.Range("Q4").Value = 3861.066777515316
$ws.Range("G5").Value = 185
$ws.Range("H5").Value = 40
$ws.Range("L5").Value = 37803.35499241169
$ws.Range("M5").Value = 9374.669557944622
$ws.Range("N5").Value = 276.8699166878031
$ws.Range("O5").Value = 366.6020013147311
$ws.Range("P5").Value = 29727.24762560007
$ws.Range("Q5").Value = 22450.99738664972
$ws.Range("G6").Value = 185
$ws.Range("H6").Value = 40
$ws.Range("L6").Value = 85268.60956677426
$ws.Range("M6").Value = 14079.44017614436
$ws.Range("N6").Value = 1712.765360158157
$ws.Range("O6").Value = 909.5129508948302
$ws.Range("P6").Value = 59190.63665068014
$ws.Range("Q6").Value = 111465.8917184739
$ws.Range("G7").Value = 146
$ws.Range("H7").Value = 24
$ws.Range("L7").Value = 46050.90522239159
$ws.Range("M7").Value = 9191.79744957916
$ws.Range("G8").Value = 146
$ws.Range("H8").Value = 24
$ws.Range("L8").Value = 65773.92710294698
$ws.Range("M8").Value = 10985.20703180743
$ws.Range("G9").Value = 146
$ws.Range("H9").Value = 24
$ws.Range("L9").Value = 127766.0304164031
$ws.Range("M9").Value = 15310.48896541679
$ws.Range("G10").Value = 146
$ws.Range("H10").Value = 24
$ws.Range("L10").Value = 83593.06812304401
$ws.Range("M10").Value = 12384.15824045096
$ws.Range("G11").Value = 146
$ws.Range("H11").Value = 24
$ws.Range("L11").Value = 99451.95742906225
$ws.Range("M11").Value = 13507.90593264003
$ws.Range("G12").Value = 130
$ws.Range("H12").Value = 70
$ws.Range("L12").Value = 1275.209927873813
$ws.Range("M12").Value = 1443.336022352533
$ws.Range("N12").Value = 20.58771781284712
$ws.Range("O12").Value = 116.4296988845964
$ws.Range("P12").Value = 10323.36646382974
$ws.Range("Q12").Value = 1825.432494217799
